$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users_LoginData")
Write-Host $ws.Name
